$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update formulas so the errors are resolved
$ws.Range("D1").Formula = "=SUM(C1,C2,D3)"
$ws.Range("D3").Formula = "=SUM(A1,76)"
$ws.Range("E2").Formula = "=SUM(A2:A5)"

# Update the active selection to E2
$ws.Range("E2").Select()
